$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 419, shifting existing rows 419-507 down to 420-508.
$ws.Rows("419").Insert()

# Populate the newly inserted row 419 with the new record's data.
$ws.Range("A419").Value = 1
$ws.Range("B419").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C419").Value = "Arica y Parinacota"
$ws.Range("D419").Value = 45209
$ws.Range("E419").Value = 15
$ws.Range("F419").Value = "Fruta"
$ws.Range("G419").Value = 100102
$ws.Range("H419").Value = "Cítricos"
$ws.Range("I419").Value = 100102003
$ws.Range("J419").Value = "Limón"
$ws.Range("K419").Value = "Sin especificar"
$ws.Range("L419").Value = "2a amarillo"
$ws.Range("M419").Value = 270
$ws.Range("N419").Value = 14000
$ws.Range("O419").Value = 15000
$ws.Range("P419").Value = 14500
$ws.Range("Q419").Value = "`$/caja 20 kilos"
$ws.Range("R419").Value = "Región de O'Higgins"
$ws.Range("S419").Value = 725
$ws.Range("T419").Value = 20
